$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44791
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 8500
$ws.Range("L2").Value = 9000
$ws.Range("M2").Value = 8750
$ws.Range("O2").Value = "Región Metropolitana"
$ws.Range("P2").Value = 875

# Row 3
$ws.Range("D3").Value = 44784
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 100
$ws.Range("K3").Value = 8000
$ws.Range("L3").Value = 9000
$ws.Range("M3").Value = 8500
$ws.Range("O3").Value = "Región Metropolitana"
$ws.Range("P3").Value = 850

# Row 4
$ws.Range("D4").Value = 44782
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 120
$ws.Range("K4").Value = 8000
$ws.Range("L4").Value = 9000
$ws.Range("M4").Value = 8500
$ws.Range("O4").Value = "Región Metropolitana"
$ws.Range("P4").Value = 850

# Row 5
$ws.Range("D5").Value = 44804
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 80
$ws.Range("K5").Value = 7000
$ws.Range("L5").Value = 7500
$ws.Range("M5").Value = 7250
$ws.Range("O5").Value = "Provincia de Diguillín"
$ws.Range("P5").Value = 725

# Row 6
$ws.Range("D6").Value = 44211
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 28
$ws.Range("K6").Value = 8000
$ws.Range("L6").Value = 8500
$ws.Range("M6").Value = 8214
$ws.Range("O6").Value = "Región Metropolitana"
$ws.Range("P6").Value = 821

# Row 7
$ws.Range("D7").Value = 44838
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 120
$ws.Range("K7").Value = 6500
$ws.Range("L7").Value = 7000
$ws.Range("M7").Value = 6750
$ws.Range("O7").Value = "Provincia de Diguillín"
$ws.Range("P7").Value = 675

# Row 8
$ws.Range("D8").Value = 44847
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 100
$ws.Range("K8").Value = 6500
$ws.Range("L8").Value = 7000
$ws.Range("M8").Value = 6750
$ws.Range("O8").Value = "Provincia de Diguillín"
$ws.Range("P8").Value = 675

# Row 9
$ws.Range("D9").Value = 44819
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 100
$ws.Range("K9").Value = 7000
$ws.Range("L9").Value = 8000
$ws.Range("M9").Value = 7500
$ws.Range("O9").Value = "Provincia de Diguillín"
$ws.Range("P9").Value = 750

# Row 10
$ws.Range("D10").Value = 44798
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 80
$ws.Range("K10").Value = 7000
$ws.Range("L10").Value = 7000
$ws.Range("M10").Value = 7000
$ws.Range("O10").Value = "Provincia de Diguillín"
$ws.Range("P10").Value = 700

# Row 11
$ws.Range("D11").Value = 44831
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 60
$ws.Range("K11").Value = 7000
$ws.Range("L11").Value = 7500
$ws.Range("M11").Value = 7250
$ws.Range("O11").Value = "Provincia de Diguillín"
$ws.Range("P11").Value = 725

# Row 12
$ws.Range("D12").Value = 44980
$ws.Range("I12").Value = "Primera"
$ws.Range("J12").Value = 60
$ws.Range("K12").Value = 7500
$ws.Range("L12").Value = 8000
$ws.Range("M12").Value = 7750
$ws.Range("O12").Value = "Provincia de Diguillín"
$ws.Range("P12").Value = 775

# Row 13
$ws.Range("D13").Value = 44817
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 60
$ws.Range("K13").Value = 7000
$ws.Range("L13").Value = 7000
$ws.Range("M13").Value = 7000
$ws.Range("O13").Value = "Provincia de Diguillín"
$ws.Range("P13").Value = 700

# Row 14
$ws.Range("D14").Value = 44817
$ws.Range("I14").Value = "Segunda"
$ws.Range("J14").Value = 60
$ws.Range("K14").Value = 8000
$ws.Range("L14").Value = 8000
$ws.Range("M14").Value = 8000
$ws.Range("O14").Value = "Provincia de Diguillín"
$ws.Range("P14").Value = 800

# Row 15
$ws.Range("D15").Value = 44841
$ws.Range("I15").Value = "Primera"
$ws.Range("J15").Value = 60
$ws.Range("K15").Value = 6500
$ws.Range("L15").Value = 7000
$ws.Range("M15").Value = 6750
$ws.Range("O15").Value = "Provincia de Diguillín"
$ws.Range("P15").Value = 675

# Row 16
$ws.Range("D16").Value = 44806
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 120
$ws.Range("K16").Value = 7000
$ws.Range("L16").Value = 7500
$ws.Range("M16").Value = 7250
$ws.Range("O16").Value = "Provincia de Diguillín"
$ws.Range("P16").Value = 725

# Row 17
$ws.Range("D17").Value = 44790
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 60
$ws.Range("K17").Value = 8500
$ws.Range("L17").Value = 9000
$ws.Range("M17").Value = 8750
$ws.Range("O17").Value = "Región Metropolitana"
$ws.Range("P17").Value = 875

# Row 18
$ws.Range("D18").Value = 44775
$ws.Range("I18").Value = "Primera"
$ws.Range("J18").Value = 60
$ws.Range("K18").Value = 8000
$ws.Range("L18").Value = 8000
$ws.Range("M18").Value = 8000
$ws.Range("O18").Value = "Región Metropolitana"
$ws.Range("P18").Value = 800

# Row 19
$ws.Range("D19").Value = 44812
$ws.Range("I19").Value = "Primera"
$ws.Range("J19").Value = 60
$ws.Range("K19").Value = 7000
$ws.Range("L19").Value = 8000
$ws.Range("M19").Value = 7500
$ws.Range("O19").Value = "Provincia de Diguillín"
$ws.Range("P19").Value = 750

# Row 20
$ws.Range("D20").Value = 44203
$ws.Range("I20").Value = "Primera"
$ws.Range("J20").Value = 27
$ws.Range("K20").Value = 7000
$ws.Range("L20").Value = 8000
$ws.Range("M20").Value = 7556
$ws.Range("O20").Value = "Región Metropolitana"
$ws.Range("P20").Value = 756

# Row 21
$ws.Range("D21").Value = 44799
$ws.Range("I21").Value = "Primera"
$ws.Range("J21").Value = 60
$ws.Range("K21").Value = 7000
$ws.Range("L21").Value = 7000
$ws.Range("M21").Value = 7000
$ws.Range("O21").Value = "Provincia de Diguillín"
$ws.Range("P21").Value = 700

# Row 22
$ws.Range("D22").Value = 44813
$ws.Range("I22").Value = "Primera"
$ws.Range("J22").Value = 120
$ws.Range("K22").Value = 7000
$ws.Range("L22").Value = 7500
$ws.Range("M22").Value = 7250
$ws.Range("O22").Value = "Provincia de Diguillín"
$ws.Range("P22").Value = 725

# Row 23
$ws.Range("D23").Value = 44810
$ws.Range("I23").Value = "Primera"
$ws.Range("J23").Value = 60
$ws.Range("K23").Value = 7000
$ws.Range("L23").Value = 8000
$ws.Range("M23").Value = 7500
$ws.Range("O23").Value = "Provincia de Diguillín"
$ws.Range("P23").Value = 750

# Row 24
$ws.Range("D24").Value = 44846
$ws.Range("I24").Value = "Primera"
$ws.Range("J24").Value = 100
$ws.Range("K24").Value = 6500
$ws.Range("L24").Value = 7000
$ws.Range("M24").Value = 6750
$ws.Range("O24").Value = "Provincia de Diguillín"
$ws.Range("P24").Value = 675
